# Update handback/handoff timestamps to reflect the latest report generation run.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 2e0da1d1... row
$wsOverview.Range("G3").Value = "2016-08-24 14:52:01"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the 2e0da1d1... row
$wsZhCn.Range("H3").Value = "2016-08-24 14:51:55"
$wsZhCn.Range("K3").Value = "2016-08-24 14:52:32"

# de-de sheet: Correspond Handoff / Handback datetimes for the 2e0da1d1... row
$wsDeDe.Range("H3").Value = "2016-08-24 14:52:01"
$wsDeDe.Range("K3").Value = "2016-08-24 14:52:40"
